# Fill in the "業務姓名" (B column) formulas for rows 18-144 on the "名單"
# sheet. These rows already have columns A/C driven by formulas that pull
# from the hidden "系統" sheet one row below (row N on 名單 <-> row N+1 on
# 系統); column B was missing the matching name-lookup formula, which is
# added here using the same IFERROR(LEFT(...),...) pattern already present
# in rows 2-17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("名單")
$ws.Activate()

for ($row = 18; $row -le 144; $row++) {
    $srcRow = $row + 1
    $formula = '=IFERROR(LEFT(系統!C' + $srcRow + ',FIND("(",系統!C' + $srcRow + ')-1),系統!C' + $srcRow + ')'
    $ws.Range("B$row").Formula = $formula
}

# Match the author's final selection / scroll position (bottom of the
# now-complete list).
[void]$ws.Range("B144").Select()
